$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, since many values look like
# plain decimal numbers (e.g. "1.002") but are actually text using "."
# as a thousands separator. This mirrors the inlineStr text cells in the
# source workbook, preventing Excel from auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.410.37"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.644.59"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "299.16"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "0.3513"
$ws.Range("E8").Value = "  -3.20%  "
$ws.Range("D9").Value = "49.77"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").Value = "0.08072"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "1.211"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("D13").Value = "22.03"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "6.360"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "7.309"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "0.00001200"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").Value = "1.640.05"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "96.57"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "0.06996"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "6.715"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "12.35"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").Value = "23.432.32"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "2.501"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "2.910"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "153.29"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "5.213"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "132.59"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "1.827.91"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "6.865"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "2.119"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").Value = "11.39"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "0.9779"
$ws.Range("E35").Value = "  -9.77%  "
$ws.Range("D36").Value = "0.02697"
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("D37").Value = "0.08742"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "0.2428"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").Value = "5.893"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "0.06807"
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("D41").Value = "12.83"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "0.6858"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "1.294"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").Value = "15.60"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "0.6325"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").Value = "2.249"
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "0.07714"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").Value = "126.96"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "1.140"
$ws.Range("E51").Value = "  -4.77%  "
